$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("client")

# The "active" column (C2) should hold the text "false" (a shared string)
# instead of the numeric 0 that was previously stored there, while
# keeping the cell's existing number format/style. A leading apostrophe
# forces Excel to store the value as text rather than reinterpreting it
# as a boolean or number.
$ws.Range("C2").Value = "'false"
